$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.009.91"
$ws.Range("E2").Value = "  +3.19%  "
$ws.Range("D3").Value = "2.952.97"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.46%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "2.950.87"
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.32%  "
$ws.Range("E11").Value = "  +6.59%  "
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("E13").Value = "  +5.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "3.442.81"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "63.021.19"
$ws.Range("E17").Value = "  +3.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "2.952.31"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "440.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("E30").Value = "  +6.69%  "
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000102"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +17.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  +2.97%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.01%  "
$ws.Range("D45").Value = "2.708.99"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0338"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "360.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.61%  "
